# Added new "Sheet3" lookup table for Variable information (manual upload).

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the last tab -----------------------------
$sheetCount = $wb.Worksheets.Count
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws3.Name = "Sheet3"

# --- Header row ----------------------------------------------------------
$ws3.Range("A1").Value = "General Name"
$ws3.Range("B1").Value = "Specific name"
$ws3.Range("C1").Value = "VariableCode"
$ws3.Range("D1").Value = "VariableID"

# --- Data rows (row 2 intentionally left blank, matching the source) ----
$ws3.Range("A3").Value = "ECRN-100 Precipitation"
$ws3.Range("B3").Value = "mm Precip"
$ws3.Range("C3").Value = "ECRN50_Precipitation"
$ws3.Range("D3").Value = 45

$ws3.Range("A4").Value = "GS3 Moisture/Temp/EC"
$ws3.Range("B4").Value = "dS/m EC Bulk"
$ws3.Range("C4").Value = "GS3_Moisture_EC"
$ws3.Range("D4").Value = 33

$ws3.Range("A5").Value = "GS3 Moisture/Temp/EC"
$ws3.Range("B5").Value = [char]0x00B0 + "C Temp"
$ws3.Range("C5").Value = "GS3_Moisture_Temp"
$ws3.Range("D5").Value = 29

$ws3.Range("A6").Value = "GS3 Moisture/Temp/EC"
$ws3.Range("B6").Value = "m" + [char]0x00B3 + "/m" + [char]0x00B3 + " VWC"
$ws3.Range("C6").Value = "GS3_Moisture_VWC"
$ws3.Range("D6").Value = 27

$ws3.Range("A7").Value = "MPS-2 Water Potential/Temp"
$ws3.Range("B7").Value = [char]0x00B0 + "C Temp"
$ws3.Range("C7").Value = "MPS2_WaterTemp"
$ws3.Range("D7").Value = 18

$ws3.Range("A8").Value = "MPS-2 Water Potential/Temp"
$ws3.Range("B8").Value = "kPa Potential"
$ws3.Range("C8").Value = "MPS2_WPot"
$ws3.Range("D8").Value = 17

$ws3.Range("A9").Value = "MPS-6 Water Potential/Temp"
$ws3.Range("B9").Value = "kPa Potential"
$ws3.Range("C9").Value = "MPS6_WaterPotential"
$ws3.Range("D9").Value = 20

$ws3.Range("A10").Value = "MPS-6 Water Potential/Temp"
$ws3.Range("B10").Value = [char]0x00B0 + "C Temp"
$ws3.Range("C10").Value = "MPS6_WaterTemp"
$ws3.Range("D10").Value = 22

$ws3.Range("A11").Value = "PYR Solar Radiation"
$ws3.Range("B11").Value = "Solar W/m" + [char]0x00B2
$ws3.Range("C11").Value = "PYR_SolarRadiation"
$ws3.Range("D11").Value = 44

$ws3.Range("A12").Value = "DS-2 Sonic Anemometer"
$ws3.Range("B12").Value = "m/s Wind"
$ws3.Range("C12").Value = "SONIC_WindSpeed"
$ws3.Range("D12").Value = 46

$ws3.Range("A13").Value = "SRS-Ni NDVI Hemispherical"
$ws3.Range("B13").Value = "800 nm"
$ws3.Range("C13").Value = "SRS_Ni_NDVI_eighthundred"
$ws3.Range("D13").Value = 37

$ws3.Range("A14").Value = "SRS-Ni NDVI Hemispherical"
$ws3.Range("B14").Value = [char]0x03B1 + " for NDVI"
$ws3.Range("C14").Value = "SRS_Ni_NDVI_Hemi_alpha"
$ws3.Range("D14").Value = 39

$ws3.Range("A15").Value = "SRS-Ni NDVI Hemispherical"
$ws3.Range("B15").Value = "630 nm"
$ws3.Range("C15").Value = "SRS_Ni_NDVI_sixthirty"
$ws3.Range("D15").Value = 35

$ws3.Range("A16").Value = "SRS-Nr NDVI Field Stop"
$ws3.Range("B16").Value = "NDVI"
$ws3.Range("C16").Value = "SRS_Nr_NDVI"
$ws3.Range("D16").Value = 43

$ws3.Range("A17").Value = "SRS-Nr NDVI Field Stop"
$ws3.Range("B17").Value = "800 nm"
$ws3.Range("C17").Value = "SRS_Nr_NDVI_eighthundred"
$ws3.Range("D17").Value = 25

$ws3.Range("A18").Value = "SRS-Nr NDVI Field Stop"
$ws3.Range("B18").Value = "630 nm"
$ws3.Range("C18").Value = "SRS_Nr_NDVI_sixthirty"
$ws3.Range("D18").Value = 23

# --- Column widths on the new sheet (approximate target layout) ---------
$ws3.Columns.Item(1).ColumnWidth = 25.59
$ws3.Range("B1:C1").EntireColumn.ColumnWidth = 17.93
$ws3.Columns.Item(4).ColumnWidth = 9.67

# --- Sheet1: selection / scroll moved, tab no longer active -------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Columns.Item(8).ColumnWidth = 22.42
$ws1.Range("B96").Select()

# --- Make the new sheet the active / selected tab ------------------------
$ws3.Activate()
$ws3.Range("B9").Select()
